$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-03 Tuesday", "2024-12-04 Wednesday"),
    @("756÷6=", "802÷8="),
    @("691÷6=", "911÷6="),
    @("894÷5=", "717÷4="),
    @("620÷4=", "367÷6="),
    @("629÷6=", "241÷4="),
    @("432÷9=", "136÷8="),
    @("610÷5=", "829÷7="),
    @("226÷5=", "228÷8="),
    @("197÷4=", "204÷7="),
    @("638÷6=", "734÷4="),
    @("597÷6=", "585÷3="),
    @("546÷8=", "471÷3="),
    @("325÷3=", "475÷8="),
    @("867÷6=", "281÷5="),
    @("705÷4=", "701÷8="),
    @("302÷5=", "866÷5="),
    @("347÷6=", "192÷3="),
    @("448÷4=", "457÷8="),
    @("931÷8=", "855÷8="),
    @("834÷4=", "217÷5="),
    @("729÷2=", "801÷6="),
    @("351÷3=", "200÷4="),
    @("549÷3=", "472÷2="),
    @("857÷3=", "809÷3="),
    @("769÷2=", "995÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
